$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.596.90'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '1.930.59'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').Value = "'326.64"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('D7').Value = "'0.4825"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').Value = "'0.4061"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = "'0.08222"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('D11').Value = "'23.88"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.17%  '
$ws.Range('D12').Value = '1.925.46'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').Value = "'6.119"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').Value = "'7.340"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.64%  '
$ws.Range('D15').Value = "'91.73"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = "'0.06874"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('D17').Value = "'1.013"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').Value = '29.608.61'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').Value = "'5.688"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').Value = "'12.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.78%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').Value = '2.163.54'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').Value = "'156.03"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').Value = "'6.413"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').Value = "'20.09"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('D30').Value = "'120.92"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').Value = "'1.015"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('D32').Value = "'0.09611"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').Value = "'5.611"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('D34').Value = "'3.558"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = "'1.392"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').Value = "'0.06389"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.92%  '
$ws.Range('D37').Value = "'0.02291"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').Value = "'10.73"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.73%  '
$ws.Range('D41').Value = "'7.930"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').Value = "'0.1850"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = "'2.479"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.56%  '
$ws.Range('D45').Value = "'1.276"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = "'12.46"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').Value = "'0.07524"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('D48').Value = "'0.5562"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = "'1.981"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').Value = "'119.16"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('D51').Value = "'2.439"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.63%  '
